# Update NATMI LR-pair (Tgfb2-Acvr1) sheet with refreshed TPM-based values.
# A new sending/receiving cluster "Resolving-Mac" is introduced, so the table
# grows from 4 clusters (3 rows each => 12 data rows) to 4x4 = 16 data rows,
# i.e. rows 2-17 (header stays on row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (A..T):
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$data = @(
    @("ECs","Tgfb2","Acvr1","ECs",3,1,2.564153333333334,7.692460000000001,0.05249149542937438,0.05249149542937438,3,1,8.165540666666667,24.496622,0.1715865889461355,0.1715865889461355,20.93769831890222,188.43928487012,0.00900683664940801,0.009006836649408014),
    @("ECs","Tgfb2","Acvr1","FAPs",3,1,2.564153333333334,7.692460000000001,0.05249149542937438,0.05249149542937438,3,1,22.011801,66.035403,0.4625449807101323,0.4625449807101323,56.44163290682001,507.9746961613801,0.02427967774082597,0.02427967774082597),
    @("ECs","Tgfb2","Acvr1","MuSCs",3,1,2.564153333333334,7.692460000000001,0.05249149542937438,0.05249149542937438,3,1,13.51552533333333,40.546576,0.2840084918355372,0.2840084918355373,34.65587933521778,311.90291401696,0.01490803045108861,0.01490803045108862),
    @("ECs","Tgfb2","Acvr1","Resolving-Mac",3,1,2.564153333333334,7.692460000000001,0.05249149542937438,0.05249149542937438,3,1,3.895588,11.686764,0.08185993850819488,0.0818599385081949,9.988884955493335,89.89996459944001,0.00429695058805178,0.004296950588051781),
    @("FAPs","Tgfb2","Acvr1","ECs",3,1,16.59481266666667,49.78443799999999,0.3397170215679993,0.3397170215679993,3,1,8.165540666666667,24.496622,0.1715865889461355,0.1715865889461355,135.5056176853818,1219.550559168436,0.05829088493779373,0.05829088493779375),
    @("FAPs","Tgfb2","Acvr1","FAPs",3,1,16.59481266666667,49.78443799999999,0.3397170215679993,0.3397170215679993,3,1,22.011801,66.035403,0.4625449807101323,0.4625449807101323,365.281714050946,3287.535426458514,0.1571344031880738,0.1571344031880738),
    @("FAPs","Tgfb2","Acvr1","MuSCs",3,1,16.59481266666667,49.78443799999999,0.3397170215679993,0.3397170215679993,3,1,13.51552533333333,40.546576,0.2840084918355372,0.2840084918355373,224.2876109982542,2018.588498984288,0.09648251894638817,0.09648251894638817),
    @("FAPs","Tgfb2","Acvr1","Resolving-Mac",3,1,16.59481266666667,49.78443799999999,0.3397170215679993,0.3397170215679993,3,1,3.895588,11.686764,0.08185993850819488,0.0818599385081949,64.64655308651467,581.818977778632,0.02780921449574354,0.02780921449574354),
    @("MuSCs","Tgfb2","Acvr1","ECs",3,1,29.664466,88.993398,0.6072695268303631,0.6072695268303631,3,1,8.165540666666667,24.496622,0.1715865889461355,0.1715865889461355,242.2264034779507,2180.037631301556,0.1041993066797557,0.1041993066797557),
    @("MuSCs","Tgfb2","Acvr1","FAPs",3,1,29.664466,88.993398,0.6072695268303631,0.6072695268303631,3,1,22.011801,66.035403,0.4625449807101323,0.4625449807101323,652.9683223632661,5876.714901269394,0.2808894715736014,0.2808894715736014),
    @("MuSCs","Tgfb2","Acvr1","MuSCs",3,1,29.664466,88.993398,0.6072695268303631,0.6072695268303631,3,1,13.51552533333333,40.546576,0.2840084918355372,0.2840084918355373,400.9308417228054,3608.377575505248,0.1724697024527717,0.1724697024527718),
    @("MuSCs","Tgfb2","Acvr1","Resolving-Mac",3,1,29.664466,88.993398,0.6072695268303631,0.6072695268303631,3,1,3.895588,11.686764,0.08185993850819488,0.0818599385081949,115.560537776008,1040.044839984072,0.04971104612423412,0.04971104612423413),
    @("Resolving-Mac","Tgfb2","Acvr1","ECs",1,0.3333333333333333,0.025497,0.076491,0.000521956172263265,0.000521956172263265,3,1,8.165540666666667,24.496622,0.1715865889461355,0.1715865889461355,0.208196790378,1.873771113402,0.00008956067917803513,0.00008956067917803516),
    @("Resolving-Mac","Tgfb2","Acvr1","FAPs",1,0.3333333333333333,0.025497,0.076491,0.000521956172263265,0.000521956172263265,3,1,22.011801,66.035403,0.4625449807101323,0.4625449807101323,0.5612348900970001,5.051114010873,0.0002414282076310464,0.0002414282076310464),
    @("Resolving-Mac","Tgfb2","Acvr1","MuSCs",1,0.3333333333333333,0.025497,0.076491,0.000521956172263265,0.000521956172263265,3,1,13.51552533333333,40.546576,0.2840084918355372,0.2840084918355373,0.3446053494240001,3.101448144816,0.0001482399852887398,0.0001482399852887398),
    @("Resolving-Mac","Tgfb2","Acvr1","Resolving-Mac",1,0.3333333333333333,0.025497,0.076491,0.000521956172263265,0.000521956172263265,3,1,3.895588,11.686764,0.08185993850819488,0.0818599385081949,0.09932580723600001,0.893932265124,0.00004272730016544365,0.00004272730016544365)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

